# Natmi following Dr Hou advice
# Updates the Bdnf-Ngfr LR-pairs sheet to add the "Neutro" and "M2" clusters
# and expand the data from a single row to the full 2x3 sending/target matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row (columns A.. T) for sheet rows 2..7.
# NOTE: every inner array literal needs a leading "," so PowerShell keeps it
# nested instead of splatting it into the outer @() list.
$rows = @(
    ,@("Neutro","Bdnf","Ngfr","M2",1,0.3333333333333333,0.08140533333333333,0.244216,0.02119526060933596,0.03145949481834748,2,0.6666666666666666,0.2576883333333334,0.7730650000000001,0.02020718570100839,0.02972883163094001,0.02097720467111111,0.18879484204,0.00042829656711412,0.0009352540246490821)
    ,@("Neutro","Bdnf","Ngfr","Neutro",1,0.3333333333333333,0.08140533333333333,0.244216,0.02119526060933596,0.03145949481834748,2,0.6666666666666666,0.241569,0.7247070000000001,0.0189431534577567,0.02786918614186859,0.019665004968,0.176985044712,0.000401505074299797,0.0008767505170216766)
    ,@("Neutro","Bdnf","Ngfr","sCs",1,0.3333333333333333,0.08140533333333333,0.244216,0.02119526060933596,0.03145949481834748,2,1,12.2530545,24.506109,0.960849660841235,0.9424019822271914,0.9974639859240001,5.984783915544,0.02036545896792205,0.02964749027667673)
    ,@("sCs","Bdnf","Ngfr","M2",2,1,3.7593275,7.518655000000001,0.9788047393906639,0.9685405051816526,2,0.6666666666666666,0.2576883333333334,0.7730650000000001,0.02020718570100839,0.02972883163094001,0.9687348379291668,5.812409027575002,0.01977888913389426,0.02879357760629093)
    ,@("sCs","Bdnf","Ngfr","Neutro",2,1,3.7593275,7.518655000000001,0.9788047393906639,0.9685405051816526,2,0.6666666666666666,0.241569,0.7247070000000001,0.0189431534577567,0.02786918614186859,0.9081369848475003,5.448821909085002,0.0185416483834569,0.02699243562484692)
    ,@("sCs","Bdnf","Ngfr","sCs",2,1,3.7593275,7.518655000000001,0.9788047393906639,0.9685405051816526,2,1,12.2530545,24.506109,0.960849660841235,0.9424019822271914,46.06324474084876,184.252978963395,0.9404842018733128,0.9127544919505147)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    $excelRow = $startRow + $i
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($excelRow, $col).Value = $rowData[$col - 1]
    }
}
